# Regenerate merged AHB files
# - rename header row labels: "<x>_old" -> "<x>_FV2210", "<x>_new" -> "<x>_FV2304"
# - wrap the data range in a native Excel Table (ListObject) named "Table1"
# - freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2210"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value2 = $val.Substring(0, $val.Length - 4) + "_FV2304"
        }
    }
}

# Turn the populated range into a native Excel table
$dataRange = $ws.Range("A1:U58")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Freeze the top header row
$null = $ws.Activate()
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
